$wb = $excel.ActiveWorkbook

# The "想去人数" (want-to-go count) values changed for three events.
# These events appear on both the "展览" sheet and the "全部类型" sheet
# (which aggregates all entries), so update both.

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 336
    $ws.Range("F3").Value = 87
    $ws.Range("F9").Value = 326
}
